$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 5 new columns before column E (so old E:I shifts to J:N)
$ws.Range("E:I").Insert()

# New year header labels for the two header rows (row 8 and row 24)
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1391/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1392/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1393/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1394/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1395/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1391/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1392/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1393/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1394/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1395/12"

# Data rows that should be filled with "-" placeholder text in the new columns
$dashRows = 10,11,12,13,14,15,16,17,18,19,26,27
foreach ($r in $dashRows) {
    $ws.Range("E${r}:I${r}").Value = "-"
}

# The total/sum row (20) gets numeric 0 instead of the "-" placeholder
$ws.Range("E20:I20").Value = 0
